# Commit: "Added CDS All studies testcase"
#
# The shared "Sample ID" query (referenced by cell B3, the SamplesTab row)
# is rewritten to drop the sample_tumor_status / sample_type columns, and the
# workbook's selection/scroll position is moved from D2 to B3.
#
# Because Excel's shared-string table drops any string that becomes
# unreferenced and appends newly-created strings at the end, simply
# overwriting B3's value with the new query text reproduces the exact
# shared-string reordering seen in the diff (the old "File Name" query
# slides from index 13 down to 12, and the rewritten "Sample ID" query
# becomes the new entry at index 13).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newSampleIdQuery = @'
SELECT
    DISTINCT (smp.sample_id) AS "Sample ID",
    sp.participant_id AS "Participant ID", 
    s.study_name AS "Study Name",
    s.phs_accession AS Accession
FROM 
    df_participant sp
JOIN 
    df_study s ON sp."study.phs_accession" = s.phs_accession
JOIN 
    df_sample smp ON smp."participant.study_participant_id" = sp.study_participant_id
JOIN
    df_diagnosis d ON d."participant.study_participant_id" = sp.study_participant_id
JOIN
    df_program p ON p.program_acronym = s."program.program_acronym"
JOIN
    df_file f1 ON f1."sample.sample_id" = smp.sample_id
JOIN
    df_genomic_info gi ON gi."file.file_id" = f1.file_id
WHERE 
   s.phs_accession = 'phs001437' AND gi.instrument_model = 'Illumina NextSeq500'
ORDER BY 
    smp.sample_id ASC
LIMIT 100;
'@

$ws.Range("B3").Value = $newSampleIdQuery

# Update the active selection / scroll position to match the saved view.
$ws.Range("B3").Select()
